$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns K (11) and L (12) added to the table
$ws.Range("K1").Value = "Switch"
$ws.Range("L1").Value = "Max"

$ws.Range("K6").Value = "Description:"
$ws.Range("L6").Value = "Description:"

$ws.Range("E7").Value = "Kill browser"
$ws.Range("F7").Value = "Go back for number of times"
$ws.Range("K7").Value = "Špatně"
$ws.Range("L7").Value = "Max. window"

# Existing cell content change
$ws.Range("B9").Value = "Must include https://"

# Move the selection / scroll position
$null = $ws.Range("B2").Select()
